# updated start script and fixed upload users error
#
# The "Users" sheet had 3 data rows; the upload bug is fixed by replacing
# the first data row's content with a correct single user record and
# removing the two erroneous extra rows that had leaked in during upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the extra (erroneous) rows first, from the bottom up, so the
# remaining row numbers don't shift while we're deleting.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# Row 2 is the single remaining user record - rewrite its fields with the
# corrected upload data. userPin / Permissions are digit strings (one has a
# leading zero) so force text formatting before assigning, otherwise Excel
# would silently coerce them to numbers and drop the leading zero.
$ws.Cells.Item(2, 1).Value = 1                                  # A2 userID

$ws.Cells.Item(2, 2).Value = "James Bond"                       # B2 username

$ws.Cells.Item(2, 3).NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = "15"                                # C2 userPin

$ws.Cells.Item(2, 4).Value = "Graham Gibson"                    # D2 supervisor
$ws.Cells.Item(2, 5).Value = "Computer Science"                 # E2 department
$ws.Cells.Item(2, 6).Value = "Queen's University"               # F2 faculty
$ws.Cells.Item(2, 7).Value = "CMC"                               # G2 institution
$ws.Cells.Item(2, 8).Value = "Academic Machine Dependent"       # H2 rateType

$ws.Cells.Item(2, 9).NumberFormat = "@"
$ws.Cells.Item(2, 9).Value = "0100100"                           # I2 Permissions
